$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 441.2
$ws.Range("I38").Value = 441.2
$ws.Range("K38").Value = 1323.6
$ws.Range("M38").Value = -951.5999999999999

$ws.Range("H39").Value = 709.4666999999999
$ws.Range("I39").Value = 137.57143
$ws.Range("J39").Value = 1209.875
$ws.Range("K39").Value = 412.71429
$ws.Range("L39").Value = 3629.625
$ws.Range("M39").Value = -116.71429
$ws.Range("N39").Value = -4221.625

$ws.Range("H74").Value = 3368.25
$ws.Range("I74").Value = 3074.8333
$ws.Range("K74").Value = 3074.8333
$ws.Range("M74").Value = -2138.8333

$ws.Range("H77").Value = 3368.25
$ws.Range("I77").Value = 3074.8333
$ws.Range("K77").Value = 15374.1665
$ws.Range("M77").Value = -10694.1665

$ws.Range("H107").Value = 894.9
$ws.Range("I107").Value = 630
$ws.Range("K107").Value = 630
$ws.Range("M107").Value = 1290

$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()

$ws.Range("H132").Value = 729.59155
$ws.Range("I132").Value = 677.7761
$ws.Range("K132").Value = 2033.3283
$ws.Range("M132").Value = 496.6716999999999

$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

$ws.Range("H135").Value = 484.21054
$ws.Range("I135").Value = 508.6111
$ws.Range("J135").Value = 45
$ws.Range("K135").Value = 4577.4999
$ws.Range("L135").Value = 405
$ws.Range("M135").Value = -2042.4999
$ws.Range("N135").Value = -5475

$ws.Range("H136").Value = 64624.75
$ws.Range("J136").Value = 64624.75
$ws.Range("L136").Value = 64624.75
$ws.Range("N136").Value = -74824.75

$ws.Range("H138").Value = 1780.6923
$ws.Range("I138").Value = 1249.0571
$ws.Range("J138").Value = 2400.9333
$ws.Range("K138").Value = 3747.1713
$ws.Range("L138").Value = 7202.7999
$ws.Range("M138").Value = 1392.8287
$ws.Range("N138").Value = -17482.7999

$ws.Range("H139").Value = 48677.668
$ws.Range("J139").Value = 48677.668
$ws.Range("L139").Value = 48677.668
$ws.Range("N139").Value = -58957.668

$ws.Range("H141").Value = 2758
$ws.Range("J141").Value = 6662
$ws.Range("L141").Value = 19986
$ws.Range("N141").Value = -30346

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2735.9453
$ws.Range("I32").Value = 2401.884
$ws.Range("K32").Value = 2401.884
$ws.Range("M32").Value = -2114.884

$ws.Range("H74").Value = 1821.5625
$ws.Range("I74").Value = 1789.4445
$ws.Range("K74").Value = 1789.4445
$ws.Range("M74").Value = -915.4445000000001

$ws.Range("H77").Value = 1821.5625
$ws.Range("I77").Value = 1789.4445
$ws.Range("K77").Value = 8947.2225
$ws.Range("M77").Value = -4579.2225

$ws.Range("H123").Value = 63993
$ws.Range("J123").Value = 63993
$ws.Range("L123").Value = 63993
$ws.Range("N123").Value = -73793

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H122").Value = 67999
$ws.Range("J122").Value = 67999
$ws.Range("L122").Value = 67999
$ws.Range("N122").Value = -77799

$ws.Range("H134").Value = 4641.0625
$ws.Range("J134").Value = 3100
$ws.Range("L134").Value = 9300
$ws.Range("N134").Value = -14370

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1919.4231
$ws.Range("I31").Value = 1678.6875
$ws.Range("K31").Value = 1678.6875
$ws.Range("M31").Value = -1383.6875

$ws.Range("H33").Value = 4976.125
$ws.Range("I33").Value = 2039.8
$ws.Range("K33").Value = 2039.8
$ws.Range("M33").Value = -1660.8

$ws.Range("H34").Value = 1919.4231
$ws.Range("I34").Value = 1678.6875
$ws.Range("K34").Value = 1678.6875
$ws.Range("M34").Value = -1476.6875

$ws.Range("H58").Value = 1451334.5
$ws.Range("I58").Value = 1813125.4
$ws.Range("K58").Value = 1813125.4
$ws.Range("M58").Value = -1812922.4

$ws.Range("H99").Value = 2137.25
$ws.Range("I99").Value = 2085.4285
$ws.Range("J99").Value = 2500
$ws.Range("K99").Value = 2085.4285
$ws.Range("L99").Value = 2500
$ws.Range("M99").Value = -587.4285
$ws.Range("N99").Value = -5496

$ws.Range("H126").Value = 2137.25
$ws.Range("I126").Value = 2085.4285
$ws.Range("J126").Value = 2500
$ws.Range("K126").Value = 6256.2855
$ws.Range("L126").Value = 7500
$ws.Range("M126").Value = -3786.2855
$ws.Range("N126").Value = -12440

$ws.Range("H134").Value = 1603.2307
$ws.Range("I134").Value = 1431.3611
$ws.Range("J134").Value = 3665.6667
$ws.Range("K134").Value = 4294.0833
$ws.Range("L134").Value = 10997.0001
$ws.Range("M134").Value = -1759.0833
$ws.Range("N134").Value = -16067.0001

$ws.Range("H136").Value = 1451334.5
$ws.Range("I136").Value = 1813125.4
$ws.Range("K136").Value = 5439376.199999999
$ws.Range("M136").Value = -5436826.199999999

$ws.Range("H138").Value = 97590
$ws.Range("J138").Value = 97590
$ws.Range("L138").Value = 97590
$ws.Range("N138").Value = -107870

$ws.Range("H139").Value = 50000
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 50000
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 50000
$ws.Range("N139").Value = -60280
$ws.Range("M139").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H123").Value = 500000000
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

$ws.Range("H131").Value = 6768410
$ws.Range("J131").Value = 12627.265
$ws.Range("L131").Value = 37881.795
$ws.Range("N131").Value = -47961.795

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2449.25
$ws.Range("I122").Value = 2399
$ws.Range("K122").Value = 7197
$ws.Range("M122").Value = -4747

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2150.6667
$ws.Range("I132").Value = 1597.8889
$ws.Range("K132").Value = 4793.6667
$ws.Range("M132").Value = -2263.6667

$ws.Range("H136").Value = 3398.875
$ws.Range("I136").Value = 1866.1333
$ws.Range("K136").Value = 5598.3999
$ws.Range("M136").Value = -3048.3999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2965.7778
$ws.Range("I132").Value = 1708.4286
$ws.Range("J132").Value = 3765.9092
$ws.Range("K132").Value = 5125.2858
$ws.Range("L132").Value = 11297.7276
$ws.Range("M132").Value = -2595.2858
$ws.Range("N132").Value = -16357.7276
